# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G ("K") values for rows 2-12 change in this revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 3
    7  = 4
    8  = 1
    9  = 1
    10 = 1
    11 = 3
    12 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
